# Apply the edits described by the commit:
#  - Clear the stale "data aggregation and summarization..." hint text out of B9:B11
#    (its shared string becomes unused and is dropped on save).
#  - Insert two new rows above the old row 35 ("Watson Discovery" / Learn-to-Rank block):
#      * a blank spacer row (picks up the formatting of the row above it)
#      * a new "Data aggregation & summarization" / "Can you tell me what happened?" row
#  - Update the selection to match the new cursor position (B9:B11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old "data aggregation and summarization..." hint cells ---
# These were the only references to that shared string, so removing the
# text here lets the now-unused string fall out of xl/sharedStrings.xml,
# and the "data collection & manipulation..." string (referenced by
# B16:B20) shifts down to take its place.
$ws.Range("B9:B11").ClearContents()

# --- Insert two new rows before the old row 35 ---
$ws.Range("A35:A36").EntireRow.Insert()

# Row 35 stays blank (just inherits formatting from the insert); set its
# height to match the new compact spacer row.
$ws.Rows("35:35").RowHeight = 21

# Row 36 holds the new "Can you tell me what happened?" entry.
$ws.Rows("36:36").RowHeight = 42
$ws.Range("C36").Value = "Data aggregation & summarization"
$ws.Range("D36").Value = "Supervised"
$ws.Range("F36").Value = "Can you tell me what happened?"
$ws.Range("E36").Value = "Focuses on telling the story of what happened for a specific event"

# --- Update the active selection to B9:B11, matching the saved view state ---
$ws.Range("B9:B11").Select() | Out-Null
